$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new rows before row 13 (shifts old rows 13-23 down to 15-25) ---
$ws.Rows.Item(13).Insert()
$ws.Rows.Item(13).Insert()

# --- Split the merged column definition (A & B sharing width/style) into its own column ---
# Touching column 1's Hidden property (no actual change) forces the engine to split the
# col min="1" max="2" entry into separate per-column entries without altering widths.
$ws.Columns.Item(1).Hidden = $false

# --- Update "Objetivos:" answer (row 10) ---
$ws.Range("B10").Value = "Possibilitar aos alunos a execução do projeto proposto e aprovado na disciplina Trabalho de Conclusão de Curso I."
$ws.Range("C10").Value = "Possibilitar aos alunos a execução do projeto proposto e aprovado na disciplina Trabalho de Conclusão de Curso I."

# --- Fill newly inserted rows 13 & 14 with the responsible teachers ---
$ws.Range("B13").Value = "1304060 - Maria das Graças de Almeida Felipe"
$ws.Range("C13").Value = "1304060 - Maria das Graças de Almeida Felipe"
$ws.Range("B14").Value = "8853480 - Tatiane da Franca Silva"
$ws.Range("C14").Value = "8853480 - Tatiane da Franca Silva"

# New rows inherited row 12's bold formatting; make B13:C14 match the normal wrapped style
# used throughout the rest of the table (vertical-top, wrap text, non-bold).
$ws.Range("B13:C14").WrapText = $true
$ws.Range("B13:C14").VerticalAlignment = -4160
$ws.Range("B13:C14").Font.Bold = $false

# Column A of the two new rows should stay blank/default (no label there).
$ws.Range("A13:A14").Style = "Normal"

# --- Update "Programa resumido:" answer (row 15, was row 13) ---
$ws.Range("B15").Value = "Desenvolvimento do trabalho de conclusão de curso, sob orientação de um professor orientador, o qual deve abordar temas relacionados à área de engenharia bioquímica."
$ws.Range("C15").Value = "Desenvolvimento do trabalho de conclusão de curso, sob orientação de um professor orientador, o qual deve abordar temas relacionados à área de engenharia bioquímica."

# --- Update "Programa:" answer (row 17, was row 15) ---
$ws.Range("B17").Value = "Elaboração de uma monografia de conclusão de curso que apresente: (1) o tema e sua importância, (2) os objetivos, (3) a revisão bibliográfica, (4) a metodologia científica (5) o desenvolvimento do projeto, (6) a análise e a discussão dos resultados, (7) as conclusões e recomendações para trabalhos futuros e (8) referências. O documento deverá atender às normas da ABNT."
$ws.Range("C17").Value = "Elaboração de uma monografia de conclusão de curso que apresente: (1) o tema e sua importância, (2) os objetivos, (3) a revisão bibliográfica, (4) a metodologia científica (5) o desenvolvimento do projeto, (6) a análise e a discussão dos resultados, (7) as conclusões e recomendações para trabalhos futuros e (8) referências. O documento deverá atender às normas da ABNT."

# --- Update "Método:" answer (row 20, was row 18) ---
$ws.Range("B20").Value = "Análise da monografia e defesa do trabalho perante uma banca de 3 examinadores, conforme Norma para Trabalho de Conclusão de Curso do curso de Engenharia Bioquímica."
$ws.Range("C20").Value = "Análise da monografia e defesa do trabalho perante uma banca de 3 examinadores, conforme Norma para Trabalho de Conclusão de Curso do curso de Engenharia Bioquímica."

# --- Update "Critério:" answer (row 21, was row 19) ---
$ws.Range("B21").Value = "A nota da disciplina será decidida pelos docentes da banca"
$ws.Range("C21").Value = "A nota da disciplina será decidida pelos docentes da banca"

# --- Update "Norma de recuperação:" answer (row 22, was row 20) ---
$ws.Range("B22").Value = "Reapresentação do trabalho modificado para nova avaliação."
$ws.Range("C22").Value = "Reapresentação do trabalho modificado para nova avaliação."

# --- Update "Bibliografia:" answer (row 23, was row 21) ---
$ws.Range("B23").Value = "Recomendada pelo Orientador."
$ws.Range("C23").Value = "Recomendada pelo Orientador."
